# data: persist cambios en clientes, historial y docs
# Adds the new cliente row (C1005 / VioletaAvila) to the "Clientes" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clientes")

$newRow = 6

# Plain text columns - straightforward assignment is fine, Excel won't
# reinterpret these as numbers/dates.
$ws.Cells.Item($newRow, 1).Value = "C1005"
$ws.Cells.Item($newRow, 2).Value = "VioletaAvila"
$ws.Cells.Item($newRow, 3).Value = "TOXQUI"
$ws.Cells.Item($newRow, 4).Value = "Martha Ortiz"

# Date-looking columns: format as Text first so Excel keeps them as the
# literal string "2025-10-09" instead of silently converting to a date
# serial, then drop back to the Normal style so no visible formatting
# change is left behind on the cell.
$ws.Cells.Item($newRow, 5).NumberFormat = "@"
$ws.Cells.Item($newRow, 5).Value = "2025-10-09"
$ws.Cells.Item($newRow, 5).Style = "Normal"

$ws.Cells.Item($newRow, 6).NumberFormat = "@"
$ws.Cells.Item($newRow, 6).Value = "2025-10-09"
$ws.Cells.Item($newRow, 6).Style = "Normal"

$ws.Cells.Item($newRow, 7).Value = "DISPERSADO"
